# Word COM-interop script implementing the InfinityMirrorPartList.docx edits
# described by the commit "Documents updated for publishing in Make share".

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, `
        $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# 1) "Mirror - 1 sq ft" -> "1/8" Mirror - 1 sq ft"
#    (operate on paragraph 2 directly - "Mirror" also appears in the
#    unrelated "Mirror film" line, so a plain text search is ambiguous).
$d.Paragraphs(2).Range.InsertBefore("1/8” ")

# 2) "Glass - 1 sq ft" -> "1/8" Glass - 1 sq ft"
#    (paragraph 3 - "Glass" also appears later in "Glass cutter").
$d.Paragraphs(3).Range.InsertBefore("1/8” ")

# 3) Remove the "Note that the glass..." paragraph entirely.
$noteRange = $d.Content
$noteRange.Find.Execute("Note that the glass and mirror film could be replaced with") | Out-Null
$noteRange.Paragraphs(1).Range.Delete()

# 4) Plywood line gets an extra clarification appended.
Replace-Text "½” Plywood – birch works nicely 2’ x 1’ is sufficient" `
    "½” Plywood – birch works nicely 2’ x 1’ is sufficient (actually 15/32”)"

# 5) Gussets line gets "or equivalent" inserted.
Replace-Text "2x4 12” for gussets" "2x4 12” or equivalent for gussets"

# 6) Insert a new "clear lacquer spray for brass" paragraph right before the
#    standalone "glue" line (not the later "spray glue (optional)" line).
$glueRange = $d.Content
$glueRange.Find.Execute("glue", $true) | Out-Null
$glueRange.Paragraphs(1).Range.InsertBefore("clear lacquer spray for brass`r")

# 7) Polishing compound: append a clarifying note. Insert right before the
#    paragraph mark (rather than doing a blanket text replace) so the
#    existing gramStart/gramEnd proofing marks around "polishing" survive.
$polishRange = $d.Content
$polishRange.Find.Execute("polishing compound", $true) | Out-Null
$polishPara = $polishRange.Paragraphs(1)
$polishEnd = $d.Range($polishPara.Range.End - 1, $polishPara.Range.End - 1)
$polishEnd.InsertBefore(" (to remove any scratches from the brasss)")

# 8) "(only 4 switches required)" text itself is unchanged in this edit.

# 9) Resistor count 4 -> 3.
Replace-Text "resistors (4)" "resistors (3)"

# 10) Add one additional blank paragraph before the existing group of blank
#     paragraphs that precede "Tools:".
$toolsRange = $d.Content
$toolsRange.Find.Execute("Tools:", $true) | Out-Null
$toolsRange.Paragraphs(1).Range.InsertParagraphBefore() | Out-Null

# 11) Update the footer date.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("2 January 2017", $true, $false, $false, $false, $false, $true, `
    $wdFindContinue, $false, "3 October 2018", $wdReplaceAll) | Out-Null
